# ControlDbp implements PUT and GET SCRIPT
#
# The DBPA web service no longer requires a local script/properties
# directory argument on the PUT/GET manager commands, so the "FROM
# scriptDir" / "TO scriptDir" / "FROM propDir" / "TO propDir" suffixes
# are dropped from the command usage text shown on the "Commands" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

$ws.Range("A2").Value = "PUT SCRIPT scriptName"
$ws.Range("A4").Value = "GET SCRIPT scriptName"
$ws.Range("A8").Value = "PUT PROPERTIES propName"
$ws.Range("A9").Value = "GET PROPERTIES propName"
